# ajout de la colonne ordre
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exigences")

# New column M: "Ordre" header (copy formatting from the adjacent "Criticité"
# header/value cells in column L so the new column matches the existing
# header/row styling), with the first data row set to 1.
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)

$ws.Range("L2").Copy()
$ws.Range("M2").PasteSpecial(-4122)

$ws.Range("M1").Value = "Ordre"
$ws.Range("M2").Value = 1

# Match the reselected/active cell from the source edit (was L2, now M2)
[void]$ws.Range("M2").Select()
